$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("K3").Value = 2.38
$ws.Range("Q3").Value = 1.7
$ws.Range("R3").Value = 2.1
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.91
$ws.Range("AC3").Value = 13
$ws.Range("AD3").Value = 8
$ws.Range("AS3").Value = 126
$ws.Range("AX3").Value = 29
$ws.Range("AY3").Value = 34
$ws.Range("AZ3").Value = 101

# Row 4 updates
$ws.Range("G4").Value = 1.95
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 4.1
$ws.Range("J4").Value = 2.63
$ws.Range("N4").Value = 9
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7
$ws.Range("W4").Value = 7
$ws.Range("X4").Value = 9
$ws.Range("Z4").Value = 17
$ws.Range("AN4").Value = 4
$ws.Range("AO4").Value = 11
$ws.Range("AQ4").Value = 41
$ws.Range("AW4").Value = 5.5
$ws.Range("AX4").Value = 21
